# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 4
    3  = 1
    4  = 6
    5  = 6
    6  = 5
    7  = 1
    8  = 2
    9  = 1
    10 = 3
    11 = 1
    12 = 2
    13 = 1
    14 = 5
    15 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
